# Update scripts with new TPM (transcript-per-million) values.
# The underlying NATMI ligand-receptor analysis was rerun with updated
# expression data; "Ligand-expressing cells" for the ECs cluster dropped
# from 2 to 1, cascading through the detection-rate / specificity /
# expression-weight columns for every ECs-source row plus the
# specificity columns of rows that depend on the ECs ligand value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.026182
$ws.Range("H2").Value = 0.078546
$ws.Range("I2").Value = 0.02060098669457318
$ws.Range("J2").Value = 0.02060098669457318
$ws.Range("M2").Value = 0.04009133333333333
$ws.Range("N2").Value = 0.120274
$ws.Range("O2").Value = 0.01033409631432067
$ws.Range("P2").Value = 0.01033409631432067
$ws.Range("Q2").Value = 0.001049671289333333
$ws.Range("R2").Value = 0.009447041604
$ws.Range("S2").Value = 0.0002128925806717579
$ws.Range("T2").Value = 0.0002128925806717579

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.026182
$ws.Range("H3").Value = 0.078546
$ws.Range("I3").Value = 0.02060098669457318
$ws.Range("J3").Value = 0.02060098669457318
$ws.Range("O3").Value = 0.578569084147867
$ws.Range("P3").Value = 0.578569084147867
$ws.Range("Q3").Value = 0.05876734046733333
$ws.Range("R3").Value = 0.528906064206
$ws.Range("S3").Value = 0.0119190940044216
$ws.Range("T3").Value = 0.0119190940044216

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.026182
$ws.Range("H4").Value = 0.078546
$ws.Range("I4").Value = 0.02060098669457318
$ws.Range("J4").Value = 0.02060098669457318
$ws.Range("O4").Value = 0.4110968195378122
$ws.Range("P4").Value = 0.4110968195378122
$ws.Range("Q4").Value = 0.04175658088333334
$ws.Range("R4").Value = 0.37580922795
$ws.Range("S4").Value = 0.008469000109479822
$ws.Range("T4").Value = 0.008469000109479822

# Row 5
$ws.Range("I5").Value = 0.9231010325934437
$ws.Range("J5").Value = 0.9231010325934434
$ws.Range("M5").Value = 0.04009133333333333
$ws.Range("N5").Value = 0.120274
$ws.Range("O5").Value = 0.01033409631432067
$ws.Range("P5").Value = 0.01033409631432067
$ws.Range("Q5").Value = 0.04703428362111112
$ws.Range("R5").Value = 0.4233085525900001
$ws.Range("S5").Value = 0.009539414978669511
$ws.Range("T5").Value = 0.009539414978669511

# Row 6
$ws.Range("I6").Value = 0.9231010325934437
$ws.Range("J6").Value = 0.9231010325934434
$ws.Range("O6").Value = 0.578569084147867
$ws.Range("P6").Value = 0.578569084147867
$ws.Range("S6").Value = 0.5340777190035391
$ws.Range("T6").Value = 0.534077719003539

# Row 7
$ws.Range("I7").Value = 0.9231010325934437
$ws.Range("J7").Value = 0.9231010325934434
$ws.Range("O7").Value = 0.4110968195378122
$ws.Range("P7").Value = 0.4110968195378122
$ws.Range("S7").Value = 0.379483898611235
$ws.Range("T7").Value = 0.3794838986112349

# Row 8
$ws.Range("I8").Value = 0.05629798071198328
$ws.Range("J8").Value = 0.05629798071198327
$ws.Range("M8").Value = 0.04009133333333333
$ws.Range("N8").Value = 0.120274
$ws.Range("O8").Value = 0.01033409631432067
$ws.Range("P8").Value = 0.01033409631432067
$ws.Range("Q8").Value = 0.002868521536222222
$ws.Range("R8").Value = 0.025816693826
$ws.Range("S8").Value = 0.0005817887549794026
$ws.Range("T8").Value = 0.0005817887549794026

# Row 9
$ws.Range("I9").Value = 0.05629798071198328
$ws.Range("J9").Value = 0.05629798071198327
$ws.Range("O9").Value = 0.578569084147867
$ws.Range("P9").Value = 0.578569084147867
$ws.Range("S9").Value = 0.03257227113990645
$ws.Range("T9").Value = 0.03257227113990644

# Row 10
$ws.Range("I10").Value = 0.05629798071198328
$ws.Range("J10").Value = 0.05629798071198327
$ws.Range("O10").Value = 0.4110968195378122
$ws.Range("P10").Value = 0.4110968195378122
$ws.Range("S10").Value = 0.02314392081709742
$ws.Range("T10").Value = 0.02314392081709742

